# Apply updated probability values to the team-specific transition matrix
# (changes to team matrices from games pulled march 7)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1898016997167139
$ws.Range("C2").Value = 0.5637393767705382
$ws.Range("J2").Value = 0.0169971671388102
$ws.Range("O2").Value = 0.0028328611898017
$ws.Range("P2").Value = 0.1444759206798867
$ws.Range("S2").Value = 0.0821529745042493
$ws.Range("B3").Value = 0.009900990099009901
$ws.Range("C3").Value = 0.009900990099009901
$ws.Range("J3").Value = 0.009900990099009901
$ws.Range("P3").Value = 0.7623762376237624
$ws.Range("S3").Value = 0.2079207920792079
$ws.Range("J4").Value = 0.01785714285714286
$ws.Range("O4").Value = 0.01785714285714286
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.2142857142857143
$ws.Range("B6").Value = 0.05607476635514019
$ws.Range("D6").Value = 0.01401869158878505
$ws.Range("F6").Value = 0.03271028037383177
$ws.Range("J6").Value = 0.2710280373831775
$ws.Range("O6").Value = 0.009345794392523364
$ws.Range("Q6").Value = 0.1962616822429906
$ws.Range("R6").Value = 0.06074766355140187
$ws.Range("S6").Value = 0.3598130841121495
$ws.Range("B7").Value = 0.1290322580645161
$ws.Range("D7").Value = 0.01612903225806452
$ws.Range("F7").Value = 0.04301075268817205
$ws.Range("J7").Value = 0.08064516129032258
$ws.Range("O7").Value = 0.02688172043010753
$ws.Range("Q7").Value = 0.1881720430107527
$ws.Range("R7").Value = 0.09677419354838709
$ws.Range("S7").Value = 0.4193548387096774
$ws.Range("B8").Value = 0.1098484848484848
$ws.Range("D8").Value = 0.02840909090909091
$ws.Range("F8").Value = 0.05871212121212122
$ws.Range("J8").Value = 0.08901515151515152
$ws.Range("O8").Value = 0.01704545454545454
$ws.Range("Q8").Value = 0.2045454545454546
$ws.Range("R8").Value = 0.07954545454545454
$ws.Range("S8").Value = 0.4128787878787879
$ws.Range("B9").Value = 0.1485714285714286
$ws.Range("D9").Value = 0.02285714285714286
$ws.Range("F9").Value = 0.02285714285714286
$ws.Range("J9").Value = 0.09142857142857143
$ws.Range("O9").Value = 0.02857142857142857
$ws.Range("Q9").Value = 0.2171428571428571
$ws.Range("R9").Value = 0.08
$ws.Range("S9").Value = 0.3885714285714286
$ws.Range("B10").Value = 0.1227380015735641
$ws.Range("D10").Value = 0.02517702596380803
$ws.Range("F10").Value = 0.07081038552321008
$ws.Range("J10").Value = 0.1030684500393391
$ws.Range("O10").Value = 0.00865460267505901
$ws.Range("Q10").Value = 0.2478363493312352
$ws.Range("R10").Value = 0.06687647521636507
$ws.Range("S10").Value = 0.3548387096774194
$ws.Range("G11").Value = 0.1374570446735395
$ws.Range("J11").Value = 0.1065292096219931
$ws.Range("K11").Value = 0.2164948453608248
$ws.Range("L11").Value = 0.5120274914089347
$ws.Range("S11").Value = 0.0274914089347079
$ws.Range("G12").Value = 0.7483443708609272
$ws.Range("J12").Value = 0.1920529801324503
$ws.Range("K12").Value = 0.006622516556291391
$ws.Range("L12").Value = 0.02649006622516556
$ws.Range("S12").Value = 0.02649006622516556
$ws.Range("G13").Value = 0.5833333333333334
$ws.Range("J13").Value = 0.3166666666666667
$ws.Range("S13").Value = 0.1
$ws.Range("G14").Value = 0.8333333333333334
$ws.Range("J14").Value = 0.1666666666666667
$ws.Range("F15").Value = 0.02314814814814815
$ws.Range("H15").Value = 0.2361111111111111
$ws.Range("I15").Value = 0.06944444444444445
$ws.Range("J15").Value = 0.3472222222222222
$ws.Range("K15").Value = 0.05092592592592592
$ws.Range("M15").Value = 0.02314814814814815
$ws.Range("N15").Value = 0.004629629629629629
$ws.Range("O15").Value = 0.05092592592592592
$ws.Range("S15").Value = 0.1944444444444444
$ws.Range("F16").Value = 0.008368200836820083
$ws.Range("H16").Value = 0.2092050209205021
$ws.Range("I16").Value = 0.09205020920502092
$ws.Range("J16").Value = 0.4267782426778243
$ws.Range("K16").Value = 0.07112970711297072
$ws.Range("M16").Value = 0.02092050209205021
$ws.Range("N16").Value = 0.004184100418410041
$ws.Range("O16").Value = 0.06276150627615062
$ws.Range("S16").Value = 0.104602510460251
$ws.Range("F17").Value = 0.01872659176029963
$ws.Range("H17").Value = 0.2134831460674157
$ws.Range("I17").Value = 0.07677902621722846
$ws.Range("J17").Value = 0.4419475655430712
$ws.Range("K17").Value = 0.08614232209737828
$ws.Range("M17").Value = 0.02059925093632959
$ws.Range("O17").Value = 0.05243445692883895
$ws.Range("S17").Value = 0.0898876404494382
$ws.Range("F18").Value = 0.02923976608187134
$ws.Range("H18").Value = 0.1637426900584795
$ws.Range("I18").Value = 0.07602339181286549
$ws.Range("J18").Value = 0.4385964912280702
$ws.Range("K18").Value = 0.1052631578947368
$ws.Range("M18").Value = 0.005847953216374269
$ws.Range("N18").Value = 0.01169590643274854
$ws.Range("O18").Value = 0.06432748538011696
$ws.Range("S18").Value = 0.1052631578947368
$ws.Range("F19").Value = 0.01767068273092369
$ws.Range("H19").Value = 0.2321285140562249
$ws.Range("I19").Value = 0.06987951807228916
$ws.Range("J19").Value = 0.3558232931726907
$ws.Range("K19").Value = 0.1044176706827309
$ws.Range("M19").Value = 0.03052208835341366
$ws.Range("N19").Value = 0.001606425702811245
$ws.Range("O19").Value = 0.07389558232931727
$ws.Range("S19").Value = 0.1140562248995984
